# Updated cryptos list on Mon Dec  4 06:08:44 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row,
# and reflects that ARBITRUM/HuobiToken swapped ranking positions
# (row 49 <-> row 50, including their Coin name / Link / Price / Volume).
#
# Note: several new Price values (column D) are plain decimal numbers
# (e.g. "229.89"). The sheet stores Price as text (to preserve values such
# as "41.504.35" which aren't valid numbers), so those assignments are
# prefixed with a leading apostrophe to force text entry, matching how
# Excel's UI keeps a numeric-looking string as text instead of silently
# converting it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.504.35'
$ws.Range('E2').Value = '  +5.25%  '
$ws.Range('D3').Value = '2.246.86'
$ws.Range('E3').Value = '  +4.01%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''229.89'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('E6').Value = '  +2.26%  '
$ws.Range('D7').Value = '''64.82'
$ws.Range('E7').Value = '  +1.06%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('E9').Value = '  +2.90%  '
$ws.Range('D10').Value = '''0.0909'
$ws.Range('E10').Value = '  +6.04%  '
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('D12').Value = '2.579.10'
$ws.Range('E12').Value = '  +3.92%  '
$ws.Range('D13').Value = '''16.12'
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').Value = '''22.30'
$ws.Range('E14').Value = '  +0.47%  '
$ws.Range('D15').Value = '''0.826'
$ws.Range('E15').Value = '  +1.39%  '
$ws.Range('D16').Value = '''5.64'
$ws.Range('E16').Value = '  +1.64%  '
$ws.Range('D17').Value = '2.246.62'
$ws.Range('E17').Value = '  +3.74%  '
$ws.Range('D18').Value = '41.408.73'
$ws.Range('D19').Value = '''73.99'
$ws.Range('E19').Value = '  +3.06%  '
$ws.Range('D20').Value = '0.0₃0920'
$ws.Range('E20').Value = '  +8.31%  '
$ws.Range('D21').Value = '''6.16'
$ws.Range('E21').Value = '  +0.98%  '
$ws.Range('D22').Value = '''252.83'
$ws.Range('E22').Value = '  +9.35%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('E24').Value = '  +1.79%  '
$ws.Range('E25').Value = '  -7.38%  '
$ws.Range('D26').Value = '''9.76'
$ws.Range('E26').Value = '  +2.73%  '
$ws.Range('D27').Value = '''172.84'
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('D28').Value = '''0.144'
$ws.Range('E28').Value = '  +3.56%  '
$ws.Range('D29').Value = '''20.47'
$ws.Range('E29').Value = '  +3.07%  '
$ws.Range('D30').Value = '''1.43'
$ws.Range('E30').Value = '  +1.86%  '
$ws.Range('D31').Value = '''2.85'
$ws.Range('E31').Value = '  +7.12%  '
$ws.Range('D32').Value = '''0.125'
$ws.Range('E32').Value = '  +2.47%  '
$ws.Range('D33').Value = '''4.76'
$ws.Range('E33').Value = '  +3.07%  '
$ws.Range('D34').Value = '''4.89'
$ws.Range('E34').Value = '  +3.07%  '
$ws.Range('D35').Value = '''7.24'
$ws.Range('E35').Value = '  +1.86%  '
$ws.Range('D36').Value = '''0.0633'
$ws.Range('E36').Value = '  +2.52%  '
$ws.Range('E37').Value = '  +7.76%  '
$ws.Range('E38').Value = '  +1.81%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '''0.000236'
$ws.Range('E40').Value = '  +49.50%  '
$ws.Range('D41').Value = '''4.80'
$ws.Range('E41').Value = '  +12.54%  '
$ws.Range('E42').Value = '  +2.32%  '
$ws.Range('D43').Value = '''8.77'
$ws.Range('E43').Value = '  +12.86%  '
$ws.Range('D44').Value = '''17.91'
$ws.Range('E44').Value = '  +0.77%  '
$ws.Range('D45').Value = '''101.59'
$ws.Range('E45').Value = '  -2.20%  '
$ws.Range('D46').Value = '''1.22'
$ws.Range('E46').Value = '  +3.40%  '
$ws.Range('D47').Value = '1.515.26'
$ws.Range('E47').Value = '  -1.61%  '
$ws.Range('D48').Value = '''0.0942'
$ws.Range('E48').Value = '  +1.79%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').Value = '''1.11'
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').Value = '''2.82'
$ws.Range('E50').Value = '  -0.48%  '
$ws.Range('D51').Value = '''51.30'
$ws.Range('E51').Value = '  +10.63%  '
